$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell value updates per diff. D/E columns hold numeric-looking text
# (prices / percentages) that must stay text, so force NumberFormat
# to "@" before assigning - otherwise Excel auto-converts them to
# numbers and strips formatting (e.g. "26.00" -> 26, "0.37%" -> 0.0037).
$textCells = @(
    "D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "D7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "D18", "E18", "D19", "E19", "D20", "E20", "E21", "D22", "E22", "E23", "E24", "D25", "E25", "D26", "E26", "D27", "D28", "E28", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "E46", "D47", "E47", "D48", "E48", "D49", "E49", "D50", "E50"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$updates = [ordered]@{
    "D2" = "246.61"
    "E2" = "0.37%"
    "D3" = "26.00"
    "E3" = "3.28%"
    "D4" = "5.093"
    "E4" = "0.93%"
    "D5" = "0.05587"
    "E5" = "-0.21%"
    "D6" = "6.478"
    "E6" = "-1.19%"
    "D7" = "0.8136"
    "D8" = "0.8467"
    "E8" = "1.39%"
    "B9" = "WazirX"
    "C9" = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
    "D9" = "0.1331"
    "E9" = "-0.30%"
    "B10" = "MandalaExchangeToken"
    "C10" = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
    "D10" = "0.06939"
    "E10" = "-0.26%"
    "B11" = "BitrueCoin"
    "C11" = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
    "D11" = "0.02817"
    "E11" = "-0.40%"
    "B12" = "BitMartToken"
    "C12" = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
    "D12" = "0.09388"
    "E12" = "-0.01%"
    "B13" = "BitForexToken"
    "C13" = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
    "D13" = "0.001511"
    "E13" = "0.10%"
    "B14" = "One"
    "C14" = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
    "D14" = "0.0005960"
    "E14" = "-93.84%"
    "B15" = "TigerCash"
    "C15" = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
    "D15" = "0.006251"
    "E15" = "1.88%"
    "B16" = "LEO"
    "C16" = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
    "D16" = "3.600"
    "E16" = "2.87%"
    "B17" = "GateToken"
    "C17" = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
    "D17" = "3.022"
    "E17" = "0.34%"
    "B18" = "BTSEToken"
    "C18" = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
    "D18" = "2.055"
    "E18" = "-1.73%"
    "B19" = "BitpandaEcosystemToken"
    "C19" = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
    "D19" = "0.3112"
    "E19" = "-2.40%"
    "D20" = "0.03168"
    "E20" = "-2.52%"
    "E21" = "-1.97%"
    "D22" = "3.748"
    "E22" = "-0.20%"
    "E23" = "-1.10%"
    "E24" = "2.50%"
    "D25" = "0.001247"
    "E25" = "0.35%"
    "D26" = "0.004551"
    "E26" = "6.11%"
    "D27" = "0.00009602"
    "D28" = "0.0001938"
    "E28" = "-0.10%"
    "D40" = "0.03654"
    "E40" = "-0.27%"
    "B41" = "KickToken"
    "C41" = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
    "D41" = "0.006151"
    "E41" = "-1.11%"
    "B42" = "BKEXToken"
    "C42" = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
    "D42" = "0.1053"
    "E42" = "-0.09%"
    "B43" = "CEJI"
    "C43" = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
    "D43" = "0.002589"
    "E43" = "-0.43%"
    "D44" = "0.007975"
    "E44" = "-3.07%"
    "D45" = "0.00005396"
    "E45" = "1.93%"
    "E46" = "0.02%"
    "D47" = "0.1450"
    "E47" = "-19.44%"
    "D48" = "0.002399"
    "E48" = "19.01%"
    "D49" = "0.00002100"
    "E49" = "0.02%"
    "D50" = "0.0002000"
    "E50" = "0.02%"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
